$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "1" to "Tsalenjikha"
$ws.Name = "Tsalenjikha"

$dots = "..."
$ellipsis = [char]0x2026

# Row 6 (Urban): B6 gets literal "...", C6:O6 get the ellipsis char (already used elsewhere)
$ws.Range("B6").Value = $dots
$ws.Range("C6:O6").Value = $ellipsis

# Row 7 (Rural): B7 gets literal "...", C7:K7 ellipsis, L7 stays 4, M7:O7 ellipsis
$ws.Range("B7").Value = $dots
$ws.Range("C7:K7").Value = $ellipsis
$ws.Range("L7").Value = 4
$ws.Range("M7:O7").Value = $ellipsis

# Delete the empty row 8 so the footnote (old row 9) shifts up to row 8
$ws.Rows.Item(8).Delete()
